# Refresh the crypto price/volume snapshot pulled in by the GitHub Actions job.
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h). All four columns hold plain
# text in the source sheet (prices use "." as a thousands separator in spots,
# e.g. "63.210.87", and some are multi-dot, so they are never numeric values).
#
# Most new Price strings still parse as a plain float (e.g. "416.49"), and Excel
# auto-converts a typed-in value that looks numeric into a Number - which would
# also silently drop significant trailing zeros ("4.50" -> 4.5). Prepending a
# leading apostrophe is the normal Excel way to force such an entry to stay text,
# exactly like the source data, without touching any other formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '63.210.87'
$ws.Range("E2").Value = '  +2.12%  '

# Row 3
$ws.Range("D3").Value = '3.489.11'
$ws.Range("E3").Value = '  +2.28%  '

# Row 4
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").Value = "'" + '416.49'
$ws.Range("E5").Value = '  +1.59%  '

# Row 6
$ws.Range("D6").Value = "'" + '131.97'
$ws.Range("E6").Value = '  +2.33%  '

# Row 7
$ws.Range("E7").Value = '  -0.99%  '

# Row 8
$ws.Range("E8").Value = '  +0.03%  '

# Row 9
$ws.Range("D9").Value = "'" + '0.737'
$ws.Range("E9").Value = '  -0.05%  '

# Row 10
$ws.Range("D10").Value = "'" + '0.154'
$ws.Range("E10").Value = '  +8.49%  '

# Row 11
$ws.Range("D11").Value = "'" + '42.84'
$ws.Range("E11").Value = '  -1.39%  '

# Row 12
$ws.Range("D12").Value = "'" + '9.88'
$ws.Range("E12").Value = '  +5.53%  '

# Row 13
$ws.Range("E13").Value = '  +3.62%  '

# Row 14
$ws.Range("D14").Value = '4.046.50'
$ws.Range("E14").Value = '  +2.37%  '

# Row 15
$ws.Range("E15").Value = '  -0.17%  '

# Row 16
$ws.Range("D16").Value = "'" + '20.68'
$ws.Range("E16").Value = '  -3.26%  '

# Row 17
$ws.Range("D17").Value = '3.484.58'
$ws.Range("E17").Value = '  +1.87%  '

# Row 18
$ws.Range("D18").Value = "'" + '12.71'
$ws.Range("E18").Value = '  +1.35%  '

# Row 19
$ws.Range("E19").Value = '  +0.17%  '

# Row 20
$ws.Range("D20").Value = '63.086.67'
$ws.Range("E20").Value = '  +1.85%  '

# Row 21
$ws.Range("D21").Value = "'" + '468.44'
$ws.Range("E21").Value = '  +4.34%  '

# Row 22
$ws.Range("D22").Value = "'" + '91.14'
$ws.Range("E22").Value = '  -0.40%  '

# Row 23
$ws.Range("E23").Value = '  +3.43%  '

# Row 24
$ws.Range("E24").Value = '  +0.88%  '

# Row 25
$ws.Range("D25").Value = "'" + '10.81'
$ws.Range("E25").Value = '  +15.73%  '

# Row 26
$ws.Range("E26").Value = '  +1.41%  '

# Row 27
$ws.Range("D27").Value = "'" + '33.75'
$ws.Range("E27").Value = '  +1.70%  '

# Row 28
$ws.Range("D28").Value = "'" + '4.80'
$ws.Range("E28").Value = '  +0.17%  '

# Row 29
$ws.Range("D29").Value = "'" + '7.60'
$ws.Range("E29").Value = '  -0.79%  '

# Row 30
$ws.Range("D30").Value = "'" + '12.20'
$ws.Range("E30").Value = '  +1.25%  '

# Row 31
$ws.Range("E31").Value = '  -3.62%  '

# Row 32
$ws.Range("E32").Value = '  -0.61%  '

# Row 33
$ws.Range("E33").Value = '  -1.12%  '

# Row 34
$ws.Range("D34").Value = "'" + '41.27'
$ws.Range("E34").Value = '  -2.98%  '

# Row 36
$ws.Range("D36").Value = "'" + '58.83'
$ws.Range("E36").Value = '  +9.41%  '

# Row 37
$ws.Range("E37").Value = '  -2.61%  '

# Row 38
$ws.Range("D38").Value = "'" + '3.09'
$ws.Range("E38").Value = '  +4.27%  '

# Row 39
$ws.Range("D39").Value = "'" + '0.999'
$ws.Range("E39").Value = '  +0.09%  '

# Row 40
$ws.Range("D40").Value = "'" + '2.77'
$ws.Range("E40").Value = '  +7.15%  '

# Row 41
$ws.Range("D41").Value = "'" + '0.136'
$ws.Range("E41").Value = '  -0.93%  '

# Row 42
$ws.Range("B42").Value = 'NEARProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D42").Value = "'" + '4.50'
$ws.Range("E42").Value = '  +3.22%  '

# Row 43
$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D43").Value = "'" + '148.10'
$ws.Range("E43").Value = '  +2.94%  '

# Row 44
$ws.Range("B44").Value = 'TheGraph'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D44").Value = "'" + '0.323'
$ws.Range("E44").Value = '  +1.29%  '

# Row 45
$ws.Range("D45").Value = "'" + '3.35'
$ws.Range("E45").Value = '  -1.44%  '

# Row 46
$ws.Range("E46").Value = '  +2.98%  '

# Row 47
$ws.Range("E47").Value = '  +37.10%  '

# Row 48
$ws.Range("D48").Value = "'" + '2.40'
$ws.Range("E48").Value = '  +11.40%  '

# Row 49
$ws.Range("D49").Value = "'" + '16.50'
$ws.Range("E49").Value = '  -1.09%  '

# Row 50
$ws.Range("D50").Value = "'" + '22.30'
$ws.Range("E50").Value = '  -1.08%  '

# Row 51
$ws.Range("E51").Value = '  -2.60%  '
